# Auto-generated edit script: update cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.790.47'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +8.46%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.951.77'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +6.85%  '

$ws.Range('E4').Value = '  -0.44%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '342.46'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.70%  '

$ws.Range('E6').Value = '  -0.36%  '

$ws.Range('E7').Value = '  +4.36%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4150'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +8.68%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '47.84'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +3.70%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.08270'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +5.72%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.037'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +8.20%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '22.75'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +8.38%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.953.75'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +6.43%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.180'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +6.01%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.405'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +5.05%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '92.07'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.86%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.002'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.37%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001060'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +3.97%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06694'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.64%  '

$ws.Range('E20').Value = '  +5.70%  '

$ws.Range('E21').Value = '  -0.34%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '29.765.85'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +8.42%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.585'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +5.73%  '

$ws.Range('E24').Value = '  +4.71%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.281'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.40%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.180.43'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +6.12%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '161.95'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.04%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '20.23'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +4.69%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.182'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +7.35%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.708'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +8.10%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '122.90'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +4.40%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.008'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +8.65%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.09653'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +3.01%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.479'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +12.81%  '

$ws.Range('E35').Value = '  +3.13%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.528'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +6.33%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.06294'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +5.74%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02320'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +6.71%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '8.493'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +4.82%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.188'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +4.00%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.6103'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +6.81%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '10.74'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +8.65%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1896'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +4.26%  '

$ws.Range('E44').Value = '  -0.31%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.398'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +35.38%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.274'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.29%  '

$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5711'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +6.31%  '

$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '12.50'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +6.56%  '

$ws.Range('E49').Value = '  +4.98%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.07365'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +7.43%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '113.52'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +3.22%  '
